# Mahipal Lomror.xlsx — scrape-update commit
# - Rename Sheet1 -> "Mahipal Lomror"
# - Insert a new leading "matchNo" column
# - Expand the single sample row into 4 full match rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Mahipal Lomror"

# Every cell in the table is stored as text in the source data (even the
# numeric-looking ones like runs/balls/sr), so force text formatting on the
# whole range *before* writing values - otherwise values such as "19" or
# "79.16" would be auto-coerced to real numbers by the COM layer.
$ws.Range("A1:M5").NumberFormat = "@"

$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $ws.Cells.Item(1, $col + 1).Value = $headers[$col]
}

$rows = @(
    @("36th","Rajasthan Royals","Mahipal Lomror","c Avesh Khan b Rabada","19","24","0","1","79.16","Delhi Capitals","Abu Dhabi","September 25","Capitals won by 33 runs"),
    @("43rd","Rajasthan Royals","Mahipal Lomror","st †Bharat b Chahal","3","4","0","0","75.00","Royal Challengers Bangalore","Dubai (DSC)","September 29","RCB won by 7 wickets (with 17 balls remaining)"),
    @("32nd","Rajasthan Royals","Mahipal Lomror","c Markram b Arshdeep Singh","43","17","2","4","252.94","Punjab Kings","Dubai (DSC)","September 21","Royals won by 2 runs"),
    @("40th","Rajasthan Royals","Mahipal Lomror","","29","28","1","1","103.57","Sunrisers Hyderabad","Dubai (DSC)","September 27","Sunrisers won by 7 wickets (with 9 balls remaining)")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowValues = $rows[$r]
    $excelRow = $r + 2
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        $ws.Cells.Item($excelRow, $col + 1).Value = $rowValues[$col]
    }
}
